# "Generate Report for Handoff"
# b.md has finished handoff processing: mark it "Ready for handoff" on the
# Overview sheet and on each language sheet (zh-cn / de-de), and record the
# newly generated handoff xliff + timestamp / error detail for that row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row 3 is b.md
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-27 22:38:16"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
# row 2 = a.md
$zhcn.Range("C2").Value = "Ready for handoff"
# row 3 = b.md
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-27 22:38:12"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8247c9ae795d26dfd56c36dc98d0d0044dc51c2f/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f05a072eabff344d89d479d27cd3db4419429ba/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
# row 2 = a.md
$dede.Range("C2").Value = "Ready for handoff"
# row 3 = b.md
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-27 22:38:16"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8247c9ae795d26dfd56c36dc98d0d0044dc51c2f/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f05a072eabff344d89d479d27cd3db4419429ba/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.17
